$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '29.961.23'
Set-TextValue "E2" '  +0.41%  '
Set-TextValue "D3" '1.893.99'
Set-TextValue "E3" '  +0.06%  '
Set-TextValue "E4" '  +0.05%  '
Set-TextValue "D5" '0.7769'
Set-TextValue "E5" '  -0.13%  '
Set-TextValue "D6" '244.00'
Set-TextValue "E6" '  -0.09%  '
Set-TextValue "D8" '0.3131'
Set-TextValue "E8" '  +0.07%  '
Set-TextValue "D9" '25.87'
Set-TextValue "E9" '  +1.76%  '
Set-TextValue "D10" '0.07250'
Set-TextValue "E10" '  +0.63%  '
Set-TextValue "D11" '0.08691'
Set-TextValue "E11" '  +7.57%  '
Set-TextValue "D12" '2.130.58'
Set-TextValue "E12" '  +12.38%  '
Set-TextValue "D13" '0.7751'
Set-TextValue "E13" '  +0.91%  '
Set-TextValue "D14" '5.423'
Set-TextValue "E14" '  -1.40%  '
Set-TextValue "D15" '94.59'
Set-TextValue "E15" '  +2.33%  '
Set-TextValue "B16" 'WrappedBTC'
Set-TextValue "C16" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D16" '30.281.06'
Set-TextValue "E16" '  +1.50%  '
Set-TextValue "B17" 'Uniswap'
Set-TextValue "C17" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D17" '6.184'
Set-TextValue "E17" '  +0.16%  '
Set-TextValue "D18" '13.92'
Set-TextValue "E18" '  -0.38%  '
Set-TextValue "D19" '2.334.32'
Set-TextValue "E19" '  +9.11%  '
Set-TextValue "D20" '246.27'
Set-TextValue "E20" '  +1.03%  '
Set-TextValue "E21" '  +1.17%  '
Set-TextValue "B22" 'Chainlink'
Set-TextValue "C22" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D22" '8.136'
Set-TextValue "E22" '  +0.06%  '
Set-TextValue "B23" 'Dai'
Set-TextValue "C23" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D23" '1.000'
Set-TextValue "E23" '  +0.05%  '
Set-TextValue "E24" '  +0.03%  '
Set-TextValue "D25" '0.1662'
Set-TextValue "E25" '  +7.11%  '
Set-TextValue "D26" '9.495'
Set-TextValue "E26" '  +1.00%  '
Set-TextValue "D27" '163.39'
Set-TextValue "E27" '  +0.53%  '
Set-TextValue "D28" '18.83'
Set-TextValue "E28" '  +0.48%  '
Set-TextValue "D29" '2.055'
Set-TextValue "E29" '  +0.37%  '
Set-TextValue "E30" '  +0.28%  '
Set-TextValue "D31" '1.543'
Set-TextValue "E31" '  -0.39%  '
Set-TextValue "D32" '4.523'
Set-TextValue "E32" '  +1.10%  '
Set-TextValue "D33" '4.131'
Set-TextValue "E33" '  +0.62%  '
Set-TextValue "D34" '0.05474'
Set-TextValue "E34" '  -0.98%  '
Set-TextValue "E35" '  -1.32%  '
Set-TextValue "D36" '0.7566'
Set-TextValue "E36" '  +1.09%  '
Set-TextValue "D37" '1.007'
Set-TextValue "E37" '  +0.62%  '
Set-TextValue "D38" '2.695'
Set-TextValue "E38" '  +2.60%  '
Set-TextValue "D39" '0.01972'
Set-TextValue "E39" '  +2.77%  '
Set-TextValue "D40" '2.787'
Set-TextValue "E40" '  +0.33%  '
Set-TextValue "D41" '0.4511'
Set-TextValue "E41" '  +2.05%  '
Set-TextValue "D42" '1.112.70'
Set-TextValue "E42" '  -2.04%  '
Set-TextValue "D43" '73.67'
Set-TextValue "E43" '  +0.10%  '
Set-TextValue "D44" '6.121'
Set-TextValue "E44" '  +3.91%  '
Set-TextValue "D45" '0.8500'
Set-TextValue "E45" '  +0.02%  '
Set-TextValue "D46" '2.217.03'
Set-TextValue "E46" '  +7.97%  '
Set-TextValue "E47" '  +0.07%  '
Set-TextValue "D48" '103.81'
Set-TextValue "E48" '  +0.10%  '
Set-TextValue "D49" '1.878'
Set-TextValue "E49" '  -0.75%  '
Set-TextValue "D50" '7.617'
Set-TextValue "E50" '  +1.96%  '
Set-TextValue "D51" '9.865'
Set-TextValue "E51" '  -0.57%  '
